$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.453262894143333
$ws.Range("C2").Value = 0.464946826018322
$ws.Range("D2").Value = 1.34141229309264
$ws.Range("E2").Value = 0.465606710193433
$ws.Range("F2").Value = 1.80132753706707

$ws.Range("B3").Value = 0.536640908522899
$ws.Range("C3").Value = 0.545583290549336
$ws.Range("D3").Value = 0.870076413594506
$ws.Range("E3").Value = 0.392518932226337
$ws.Range("F3").Value = 0.374672849089424

$ws.Range("B4").Value = 0.488556234623904
$ws.Range("C4").Value = 0.597115893972384
$ws.Range("E4").Value = 0.243154295480897

$ws.Range("B5").Value = 1.62749589944488
$ws.Range("C5").Value = 1.77114502143414

$ws.Range("B6").Value = 0.522569588062163
$ws.Range("C6").Value = 0.534021610917451
$ws.Range("D6").Value = 4.52652963238878
$ws.Range("E6").Value = 0.534752274100427
$ws.Range("F6").Value = 4.39407843884091
